$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = '1050954263'
$ws.Range("D16").Value = 'JORGE LEONARDO ELLES MERCADO'
$ws.Range("E16").Value = '2112'
$ws.Range("F16").Value = 50880
$ws.Range("G16").Value = 1272000

$ws.Range("C17").Value = '1050954263'
$ws.Range("D17").Value = 'JORGE LEONARDO ELLES MERCADO'
$ws.Range("E17").Value = '2111'
$ws.Range("F17").Value = 50880
$ws.Range("G17").Value = 1272000

$ws.Range("C18").Value = '1050954263'
$ws.Range("D18").Value = 'JORGE LEONARDO ELLES MERCADO'
$ws.Range("E18").Value = '2110'
$ws.Range("F18").Value = 50880
$ws.Range("G18").Value = 1272000

$ws.Range("C19").Value = '73181762'
$ws.Range("D19").Value = 'JHOFRE BATISTA FUENTES'
$ws.Range("E19").Value = '2112'
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908526

$ws.Range("C20").Value = '73181762'
$ws.Range("D20").Value = 'JHOFRE BATISTA FUENTES'
$ws.Range("E20").Value = '2111'
$ws.Range("F20").Value = 36341
$ws.Range("G20").Value = 908526

$ws.Range("C21").Value = '73181762'
$ws.Range("D21").Value = 'JHOFRE BATISTA FUENTES'
$ws.Range("E21").Value = '2110'
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 908526

$ws.Range("C22").Value = '1047404169'
$ws.Range("D22").Value = 'RAFAEL EDUARDO PEREZ FLOREZ'
$ws.Range("E22").Value = '2112'
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 690000

$ws.Range("C23").Value = '1047404169'
$ws.Range("D23").Value = 'RAFAEL EDUARDO PEREZ FLOREZ'
$ws.Range("E23").Value = '2111'
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 690000

$ws.Range("C24").Value = '1047404169'
$ws.Range("D24").Value = 'RAFAEL EDUARDO PEREZ FLOREZ'
$ws.Range("E24").Value = '2110'
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 690000

$ws.Range("C25").Value = '1047407998'
$ws.Range("D25").Value = 'LUIS ALFONSO BELEÑO BRAVO'
$ws.Range("E25").Value = '2112'
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 908526

$ws.Range("C26").Value = '1047407998'
$ws.Range("D26").Value = 'LUIS ALFONSO BELEÑO BRAVO'
$ws.Range("E26").Value = '2111'
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 908526

$ws.Range("C27").Value = '1047407998'
$ws.Range("D27").Value = 'LUIS ALFONSO BELEÑO BRAVO'
$ws.Range("E27").Value = '2110'
$ws.Range("F27").Value = 36341
$ws.Range("G27").Value = 908526

$ws.Range("C28").Value = '1047371389'
$ws.Range("D28").Value = 'YONATAN RODRIGUEZ LUNA'
$ws.Range("E28").Value = '2112'
$ws.Range("F28").Value = 36341
$ws.Range("G28").Value = 908526

$ws.Range("C29").Value = '1047371389'
$ws.Range("D29").Value = 'YONATAN RODRIGUEZ LUNA'
$ws.Range("E29").Value = '2111'
$ws.Range("F29").Value = 36341
$ws.Range("G29").Value = 908526

$ws.Range("C30").Value = '1047371389'
$ws.Range("D30").Value = 'YONATAN RODRIGUEZ LUNA'
$ws.Range("E30").Value = '2110'
$ws.Range("F30").Value = 36341
$ws.Range("G30").Value = 908526

$ws.Range("C31").Value = '1126256862'
$ws.Range("D31").Value = 'GUSTAVO ENRIQUE BAPTISTA GUERRERO'
$ws.Range("E31").Value = '2112'
$ws.Range("F31").Value = 36341
$ws.Range("G31").Value = 908526

$ws.Range("C32").Value = '1126256862'
$ws.Range("D32").Value = 'GUSTAVO ENRIQUE BAPTISTA GUERRERO'
$ws.Range("E32").Value = '2111'
$ws.Range("F32").Value = 36341
$ws.Range("G32").Value = 908526

$ws.Range("C33").Value = '1126256862'
$ws.Range("D33").Value = 'GUSTAVO ENRIQUE BAPTISTA GUERRERO'
$ws.Range("E33").Value = '2110'
$ws.Range("F33").Value = 36341
$ws.Range("G33").Value = 908526

$ws.Range("C34").Value = '1020786453'
$ws.Range("D34").Value = 'RONALD ABRAHAM QUINTANA MONTALVO'
$ws.Range("E34").Value = '2112'
$ws.Range("F34").Value = 36341
$ws.Range("G34").Value = 908526

$ws.Range("C35").Value = '1020786453'
$ws.Range("D35").Value = 'RONALD ABRAHAM QUINTANA MONTALVO'
$ws.Range("E35").Value = '2111'
$ws.Range("F35").Value = 36341
$ws.Range("G35").Value = 908526

$ws.Range("C36").Value = '1020786453'
$ws.Range("D36").Value = 'RONALD ABRAHAM QUINTANA MONTALVO'
$ws.Range("E36").Value = '2110'
$ws.Range("F36").Value = 36341
$ws.Range("G36").Value = 908526

$ws.Range("C37").Value = '1047470076'
$ws.Range("D37").Value = 'DEIMER ANTONIO MARRUGO HERRERA'
$ws.Range("E37").Value = '2112'
$ws.Range("F37").Value = 36341
$ws.Range("G37").Value = 828116

$ws.Range("C38").Value = '1047470076'
$ws.Range("D38").Value = 'DEIMER ANTONIO MARRUGO HERRERA'
$ws.Range("E38").Value = '2111'
$ws.Range("F38").Value = 36341
$ws.Range("G38").Value = 828116

$ws.Range("C39").Value = '1047470076'
$ws.Range("D39").Value = 'DEIMER ANTONIO MARRUGO HERRERA'
$ws.Range("E39").Value = '2110'
$ws.Range("F39").Value = 29073
$ws.Range("G39").Value = 828116

$ws.Range("C40").Value = '1007470327'
$ws.Range("D40").Value = 'JULIAN SIERRA GUERRERO'
$ws.Range("E40").Value = '2112'
$ws.Range("F40").Value = 36341
$ws.Range("G40").Value = 908526

$ws.Range("C41").Value = '1007470327'
$ws.Range("D41").Value = 'JULIAN SIERRA GUERRERO'
$ws.Range("E41").Value = '2111'
$ws.Range("F41").Value = 36341
$ws.Range("G41").Value = 908526

$ws.Range("C42").Value = '1007470327'
$ws.Range("D42").Value = 'JULIAN SIERRA GUERRERO'
$ws.Range("E42").Value = '2110'
$ws.Range("F42").Value = 36341
$ws.Range("G42").Value = 908526

$ws.Range("C43").Value = '1048932784'
$ws.Range("D43").Value = 'RUBEN DARIO GUERRERO PEREZ'
$ws.Range("E43").Value = '2112'
$ws.Range("F43").Value = 36341
$ws.Range("G43").Value = 908526

$ws.Range("C44").Value = '1048932784'
$ws.Range("D44").Value = 'RUBEN DARIO GUERRERO PEREZ'
$ws.Range("E44").Value = '2111'
$ws.Range("F44").Value = 36341
$ws.Range("G44").Value = 908526

$ws.Range("C45").Value = '1048932784'
$ws.Range("D45").Value = 'RUBEN DARIO GUERRERO PEREZ'
$ws.Range("E45").Value = '2110'
$ws.Range("F45").Value = 36341
$ws.Range("G45").Value = 908526

$ws.Range("C46").Value = '1235046173'
$ws.Range("D46").Value = 'JOSE ALFREDO RIVERA ROMERO'
$ws.Range("E46").Value = '2112'
$ws.Range("F46").Value = 36341
$ws.Range("G46").Value = 908526

$ws.Range("C47").Value = '1235046173'
$ws.Range("D47").Value = 'JOSE ALFREDO RIVERA ROMERO'
$ws.Range("E47").Value = '2111'
$ws.Range("F47").Value = 36341
$ws.Range("G47").Value = 908526

$ws.Range("C48").Value = '1235046173'
$ws.Range("D48").Value = 'JOSE ALFREDO RIVERA ROMERO'
$ws.Range("E48").Value = '2110'
$ws.Range("F48").Value = 36341
$ws.Range("G48").Value = 908526

$ws.Range("C49").Value = '1047425628'
$ws.Range("D49").Value = 'KLEIBER BENJAMIN BARRIOS MARIMON'
$ws.Range("E49").Value = '2112'
$ws.Range("F49").Value = 36341
$ws.Range("G49").Value = 908526

$ws.Range("C50").Value = '1047425628'
$ws.Range("D50").Value = 'KLEIBER BENJAMIN BARRIOS MARIMON'
$ws.Range("E50").Value = '2111'
$ws.Range("F50").Value = 36341
$ws.Range("G50").Value = 908526

$ws.Range("C51").Value = '1047425628'
$ws.Range("D51").Value = 'KLEIBER BENJAMIN BARRIOS MARIMON'
$ws.Range("E51").Value = '2110'
$ws.Range("F51").Value = 36341
$ws.Range("G51").Value = 908526

$ws.Range("C52").Value = '1193556287'
$ws.Range("D52").Value = 'LUIS ANGEL GUERRERO PEREZ'
$ws.Range("E52").Value = '2112'
$ws.Range("F52").Value = 36341
$ws.Range("G52").Value = 908526

$ws.Range("C53").Value = '1193556287'
$ws.Range("D53").Value = 'LUIS ANGEL GUERRERO PEREZ'
$ws.Range("E53").Value = '2111'
$ws.Range("F53").Value = 36341
$ws.Range("G53").Value = 908526

$ws.Range("C54").Value = '1193556287'
$ws.Range("D54").Value = 'LUIS ANGEL GUERRERO PEREZ'
$ws.Range("E54").Value = '2110'
$ws.Range("F54").Value = 29073
$ws.Range("G54").Value = 908526
